$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 7335
$ws.Range("B2").Value = 107
$ws.Range("C2").Value = 280
